$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write updated cell values. NumberFormat is set to Text ("@") immediately
# before each write so values that look numeric (e.g. "1.00", "9.19") are
# preserved verbatim as text, matching the original inline-string cells,
# instead of Excel silently re-interpreting them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.113.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.320.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.76"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.348.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.045.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.21%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.25"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.15"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.97"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.90"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.993.15"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.52"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +8.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.05"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.547.77"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.50%  "
